$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 12 (format + values) into the new row 13, then overwrite the
# cells that actually changed for this run (new ANN hyper-parameters /
# scores), matching the appended results row in the updated workbook.
$ws.Range("A12:J12").Copy($ws.Range("A13:J13"))

$ws.Range("F13").Value = "Layer=100*tanh+ 100*tanh  + 100 * sigmoid + softmax`nlearning_rate=0,1`nn_stable=2`nn_iter=21"
$ws.Range("G13").Value = 0.23028
$ws.Range("H13").Value = 2.55942
$ws.Range("I13").Value = 2.56269
$ws.Range("J13").Value = "505/1390"

# Row 12 had an auto height of 49.25 for its 4-line wrapped note; row 13's
# note wraps to more lines, so it grows accordingly.
$ws.Rows.Item(13).RowHeight = 61.15

# Move the active selection on to the next empty row, as Excel does after
# entering data in the last cell of a row.
$ws.Range("J14").Select() | Out-Null
